$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 19773.627
$ws.Range("J87").Value = 19773.627
$ws.Range("L87").Value = 19773.627
$ws.Range("N87").Value = -22269.627
$ws.Range("H90").Value = 19773.627
$ws.Range("J90").Value = 19773.627
$ws.Range("L90").Value = 59320.881
$ws.Range("N90").Value = -71800.88099999999
$ws.Range("H132").Value = 22419.215
$ws.Range("I132").Value = 23759.154
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 71277.462
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -68747.462
$ws.Range("N132").Value = -20060
$ws.Range("H135").Value = 1166.25
$ws.Range("I135").Value = 613.4286
$ws.Range("J135").Value = 5036
$ws.Range("K135").Value = 5520.8574
$ws.Range("L135").Value = 45324
$ws.Range("M135").Value = -2985.8574
$ws.Range("N135").Value = -50394

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5176.8125
$ws.Range("I32").Value = 3121.9333
$ws.Range("K32").Value = 3121.9333
$ws.Range("M32").Value = -2834.9333
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H132").Value = 3252.8462
$ws.Range("I132").Value = 1594
$ws.Range("J132").Value = 4469.3335
$ws.Range("K132").Value = 4782
$ws.Range("L132").Value = 13408.0005
$ws.Range("M132").Value = -2252
$ws.Range("N132").Value = -18468.0005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 3812.261
$ws.Range("I80").Value = 2615.9092
$ws.Range("J80").Value = 4908.9165
$ws.Range("K80").Value = 2615.9092
$ws.Range("L80").Value = 4908.9165
$ws.Range("M80").Value = -1617.9092
$ws.Range("N80").Value = -6904.9165
$ws.Range("H83").Value = 3812.261
$ws.Range("I83").Value = 2615.9092
$ws.Range("J83").Value = 4908.9165
$ws.Range("K83").Value = 13079.546
$ws.Range("L83").Value = 24544.5825
$ws.Range("M83").Value = -8087.546
$ws.Range("N83").Value = -34528.5825

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3941.7144
$ws.Range("I31").Value = 1516.5294
$ws.Range("K31").Value = 1516.5294
$ws.Range("M31").Value = -1221.5294
$ws.Range("H34").Value = 3941.7144
$ws.Range("I34").Value = 1516.5294
$ws.Range("K34").Value = 1516.5294
$ws.Range("M34").Value = -1314.5294
$ws.Range("H58").Value = 2327.56
$ws.Range("I58").Value = 1469.6428
$ws.Range("J58").Value = 3419.4546
$ws.Range("K58").Value = 1469.6428
$ws.Range("L58").Value = 3419.4546
$ws.Range("M58").Value = -1266.6428
$ws.Range("N58").Value = -3825.4546
$ws.Range("H134").Value = 2351.4614
$ws.Range("I134").Value = 2187.0908
$ws.Range("J134").Value = 3255.5
$ws.Range("K134").Value = 6561.2724
$ws.Range("L134").Value = 9766.5
$ws.Range("M134").Value = -4026.2724
$ws.Range("N134").Value = -14836.5
$ws.Range("H136").Value = 2327.56
$ws.Range("I136").Value = 1469.6428
$ws.Range("J136").Value = 3419.4546
$ws.Range("K136").Value = 4408.928400000001
$ws.Range("L136").Value = 10258.3638
$ws.Range("M136").Value = -1858.928400000001
$ws.Range("N136").Value = -15358.3638

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 7217.5
$ws.Range("I106").Value = 4435
$ws.Range("J106").Value = 10000
$ws.Range("K106").Value = 13305
$ws.Range("L106").Value = 30000
$ws.Range("M106").Value = -12359
$ws.Range("N106").Value = -31892
$ws.Range("H131").Value = 2925.6365
$ws.Range("J131").Value = 3484.6667
$ws.Range("L131").Value = 10454.0001
$ws.Range("N131").Value = -20534.0001
$ws.Range("H140").Value = 1448.3334
$ws.Range("I140").Value = 1448.3334
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 4345.0002
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 834.9997999999996
$ws.Range("N140").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 33390.656
$ws.Range("I122").Value = 1889.1364
$ws.Range("K122").Value = 5667.4092
$ws.Range("M122").Value = -3217.4092
$ws.Range("H134").Value = 15481.5
$ws.Range("J134").Value = 15481.5
$ws.Range("L134").Value = 46444.5
$ws.Range("N134").Value = -51514.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 59816.668
$ws.Range("I40").Value = 86750
$ws.Range("J40").Value = 5950
$ws.Range("K40").Value = 86750
$ws.Range("L40").Value = 5950
$ws.Range("M40").Value = -86614
$ws.Range("N40").Value = -6222
$ws.Range("H122").Value = 4057
$ws.Range("I122").Value = 2500.8
$ws.Range("J122").Value = 5613.2
$ws.Range("K122").Value = 7502.400000000001
$ws.Range("L122").Value = 16839.6
$ws.Range("M122").Value = -5052.400000000001
$ws.Range("N122").Value = -21739.6
$ws.Range("H136").Value = 3314.9443
$ws.Range("I136").Value = 2261.7144
$ws.Range("J136").Value = 7001.25
$ws.Range("K136").Value = 6785.1432
$ws.Range("L136").Value = 21003.75
$ws.Range("M136").Value = -4235.1432
$ws.Range("N136").Value = -26103.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 32495.818
$ws.Range("J54").Value = 32495.818
$ws.Range("L54").Value = 32495.818
$ws.Range("N54").Value = -33535.818
$ws.Range("H62").Value = 4850.3
$ws.Range("I62").Value = 4633.3335
$ws.Range("J62").Value = 4943.2856
$ws.Range("K62").Value = 4633.3335
$ws.Range("L62").Value = 4943.2856
$ws.Range("M62").Value = -4009.3335
$ws.Range("N62").Value = -6191.2856
$ws.Range("H65").Value = 4850.3
$ws.Range("I65").Value = 4633.3335
$ws.Range("J65").Value = 4943.2856
$ws.Range("K65").Value = 23166.6675
$ws.Range("L65").Value = 24716.428
$ws.Range("M65").Value = -20046.6675
$ws.Range("N65").Value = -30956.428
$ws.Range("H81").Value = 888.8889
$ws.Range("I81").Value = 857.1429000000001
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 1714.2858
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -653.2858000000001
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 888.8889
$ws.Range("I84").Value = 857.1429000000001
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 8571.429
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -3267.429
$ws.Range("N84").Value = -20608
$ws.Range("H132").Value = 2211.963
$ws.Range("I132").Value = 1957.3125
$ws.Range("J132").Value = 2582.3635
$ws.Range("K132").Value = 5871.9375
$ws.Range("L132").Value = 7747.0905
$ws.Range("M132").Value = -3341.9375
$ws.Range("N132").Value = -12807.0905
